# This workbook's data rows (4, 5, 6) are rotated: the record that was in
# row 5 moves to row 4, the record in row 6 moves to row 5, and the record
# that was in row 4 moves to row 6 (a cyclic "shift up" of the three
# observation records, columns A,B,D,E,F,G,H,I,J,Q,R).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Capture the "before" values for the columns that change, for rows 4-6.
$cols = @("A","B","D","E","F","G","H","I","J","Q","R")

$row4 = @{}
$row5 = @{}
$row6 = @{}

foreach ($col in $cols) {
    $row4[$col] = $ws.Range($col + "4").Value2
    $row5[$col] = $ws.Range($col + "5").Value2
    $row6[$col] = $ws.Range($col + "6").Value2
}

# "Antal" (column I) is stored as text in the source data, not a number,
# so make sure the destination cell is formatted as text before writing it.
$ws.Range("I4:I6").NumberFormat = "@"

# Apply the rotation: new row4 = old row5, new row5 = old row6, new row6 = old row4
foreach ($col in $cols) {
    $ws.Range($col + "4").Value = $row5[$col]
    $ws.Range($col + "5").Value = $row6[$col]
    $ws.Range($col + "6").Value = $row4[$col]
}

# I4/J4 and I6/J6 end up empty (no "Antal"/"Enhet" for those records) -
# clear them explicitly in case Value2 returned $null and didn't blank the cell.
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("I6").ClearContents()
$ws.Range("J6").ClearContents()
